# This edit inserts one new data row into the weekly price table.
# The new row is inserted at row 123 (pushing the existing rows 123..199
# down to 124..200) and is populated with a new "Zapallo italiano" price
# observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 123; this shifts all the
# rows below (old 123..199) down by one (new 124..200), carrying their
# formatting (incl. the date style on column D) along with them.
$ws.Rows.Item(123).Insert()

# Populate the newly inserted row 123 with the new data point.
$ws.Range("A123").Value2 = 6
$ws.Range("B123").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C123").Value2 = "Metropolitana"
$ws.Range("D123").Value2 = 44438
$ws.Range("E123").Value2 = 13
$ws.Range("F123").Value2 = 100112032
$ws.Range("G123").Value2 = "Zapallo italiano"
$ws.Range("H123").Value2 = "Sin especificar"
$ws.Range("I123").Value2 = "Primera"
$ws.Range("J123").Value2 = 250
$ws.Range("K123").Value2 = 13000
$ws.Range("L123").Value2 = 14000
$ws.Range("M123").Value2 = 13400
$ws.Range("N123").Value2 = "`$/caja 50 unidades"
$ws.Range("O123").Value2 = "Región de Arica y Parinacota"
$ws.Range("P123").Value2 = 268
$ws.Range("Q123").Value2 = 50
$ws.Range("R123").Value2 = "Hortaliza"
